# Trade #36 closed at 2026-02-17 12:39:04 - unknown UNKNOWN +0.000%
# Applies the new closed trade to the "All Trades" and "MarketMaking" sheets,
# and refreshes the aggregate stats on "Summary" and "Strategy Status".

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$strategyStatus = $wb.Worksheets.Item("Strategy Status")
$allTrades = $wb.Worksheets.Item("All Trades")
$marketMaking = $wb.Worksheets.Item("MarketMaking")

# --- Summary sheet updates ---
$summary.Range("B5").Value = 0.41    # Total P&L %
$summary.Range("B6").Value = 36      # Total Trades
$summary.Range("B9").Value = 38.89   # Win Rate %

# --- Strategy Status sheet updates (MarketMaking row) ---
$strategyStatus.Range("D4").Value = 36     # Trades
$strategyStatus.Range("G4").Value = 38.89  # Win Rate %

# --- New trade row data (row 37 on both trade-log sheets) ---
# Note: the date string is prefixed with a leading apostrophe so Excel keeps
# it as literal text (matching the existing Date column cells) instead of
# auto-converting it to a date serial number.
$tradeRow = @(36, "'2026-02-17", "12:38:57", "MarketMaking", "DOWN", 0.07000000000000001, 0.07169200000000001, "CLOSED", 2.417, 0, 100.74, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

foreach ($sheet in @($allTrades, $marketMaking)) {
    $sheet.Cells.Item(37, 1).Value = $tradeRow[0]
    $sheet.Cells.Item(37, 2).Value = $tradeRow[1]
    $sheet.Cells.Item(37, 3).Value = $tradeRow[2]
    $sheet.Cells.Item(37, 4).Value = $tradeRow[3]
    $sheet.Cells.Item(37, 5).Value = $tradeRow[4]
    $sheet.Cells.Item(37, 6).Value = $tradeRow[5]
    $sheet.Cells.Item(37, 7).Value = $tradeRow[6]
    $sheet.Cells.Item(37, 8).Value = $tradeRow[7]
    $sheet.Cells.Item(37, 9).Value = $tradeRow[8]
    $sheet.Cells.Item(37, 10).Value = $tradeRow[9]
    $sheet.Cells.Item(37, 11).Value = $tradeRow[10]
    $sheet.Cells.Item(37, 12).Value = $tradeRow[11]
    $sheet.Cells.Item(37, 13).Value = $tradeRow[12]
    $sheet.Cells.Item(37, 14).Value = $tradeRow[13]
    $sheet.Cells.Item(37, 15).Value = $tradeRow[14]
    $sheet.Cells.Item(37, 16).Value = $tradeRow[15]
    $sheet.Cells.Item(37, 17).Value = $tradeRow[16]
}
